{"js": "// The author removed a stray extra blank \"Bold heading\" paragraph from the\n// run of empty paragraphs near the end of the document (the \"wonky font\n// formatting\" fix). Word automatically keeps a single `_GoBack` bookmark\n// around the most recent edit location, so deleting that paragraph makes\n// Word relocate the bookmark from its old spot (mid-sentence, after\n// \"...assign you a grader for \") to the end of the paragraph that is now\n// last in that empty-paragraph run (immediately before the closing\n// \"Normal (Web)\" copyright paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"style,text\");\n}\nawait context.sync();\n\n// Locate the run of consecutive, empty \"Bold heading\" paragraphs that sits\n// right before the very last paragraph of the document (the copyright\n// notice). That run is the set of blank spacer paragraphs at the end of\n// the handout.\nconst items = paragraphs.items;\nconst lastIndex = items.length - 1;\n\nlet runEnd = -1;\nfor (let i = lastIndex - 1; i >= 0; i--) {\n  if (items[i].style === \"Bold heading\" && items[i].text === \"\") {\n    runEnd = i;\n  } else {\n    break;\n  }\n}\n\nif (runEnd === -1) {\n  throw new Error(\"Could not find the trailing blank 'Bold heading' paragraphs.\");\n}\n\n// Paragraph just before the one we delete keeps going forward - it becomes\n// the new last paragraph in the run, and receives the relocated bookmark.\nconst deleteIndex = lastIndex - 1;\nconst keepIndex = deleteIndex - 1;\n\n// Remove the old `_GoBack` bookmark from wherever it currently lives in the\n// document (mid-paragraph, right after \"...a grader for \"). `deleteBookmark`\n// is a no-op if the bookmark isn't present, so this is safe either way.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-create the bookmark at the end of the paragraph that will become the\n// new last paragraph of the blank run.\nitems[keepIndex].getRange(\"End\").insertBookmark(\"_GoBack\");\n\n// Delete the extra blank paragraph.\nitems[deleteIndex].delete();\n\nawait context.sync();\n", "ps1": "# The author removed a stray extra blank \"Bold heading\" paragraph from the\n# run of empty paragraphs near the end of the document (the \"wonky font\n# formatting\" fix). Word automatically keeps a single `_GoBack` bookmark\n# around the most recent edit location, so deleting that paragraph makes\n# Word relocate the bookmark from its old spot (mid-sentence, after\n# \"...assign you a grader for \") to the end of the paragraph that is now\n# last in that empty-paragraph run (immediately before the closing\n# \"Normal (Web)\" copyright paragraph).\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# Walk backward (skipping the very last paragraph, the copyright notice)\n# while we're inside the trailing run of blank \"Bold heading\" spacer\n# paragraphs, to find where that run starts.\n$runStart = -1\nfor ($i = $count - 1; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $isBlankHeading = ($p.Style.NameLocal -eq \"Bold heading\") -and ($p.Range.Text -eq [char]13)\n    if ($isBlankHeading) {\n        $runStart = $i\n    } else {\n        break\n    }\n}\n\nif ($runStart -eq -1) {\n    throw \"Could not find the trailing blank 'Bold heading' paragraphs.\"\n}\n\n# The last paragraph of that run is the extra one to remove; the one right\n# before it becomes the new last paragraph of the run and receives the\n# relocated bookmark.\n$deleteIndex = $count - 1\n$keepIndex = $deleteIndex - 1\n\n# Remove the old `_GoBack` bookmark from wherever it currently lives in the\n# document (mid-paragraph, right after \"...a grader for \").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Re-create the bookmark at the end of the paragraph that will become the\n# new last paragraph of the blank run (collapsed, right before its\n# paragraph mark).\n$keepRange = $d.Paragraphs.Item($keepIndex).Range\n$keepRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $keepRange)\n\n# Delete the extra blank paragraph.\n$d.Paragraphs.Item($deleteIndex).Range.Delete()\n"}
